$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the RF (column I) values for rows 37 through 69
# from 82.6059574468085 to 103.26125
$ws.Range("I37:I69").Value = 103.26125
